$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Background Processing" flag for row 2 to TRUE
$ws.Range("G2").Value = $true

# Delete row 3 entirely (it was a duplicate of row 2 except G3=TRUE)
$ws.Rows("3").Delete()

# Select the whole of row 2 as the active selection
$ws.Rows("2").Select()
